$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 23.496322
$ws.Range("H2").Value = 70.488966
$ws.Range("I2").Value = 0.1321353991144917
$ws.Range("J2").Value = 0.1321353991144917
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8369776666666665
$ws.Range("N2").Value = 2.510933
$ws.Range("O2").Value = 0.0694586718035551
$ws.Range("P2").Value = 0.06945867180355511
$ws.Range("Q2").Value = 19.66589676280866
$ws.Range("R2").Value = 176.993070865278
$ws.Range("S2").Value = 0.009177949320725241
$ws.Range("T2").Value = 0.009177949320725243

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 23.496322
$ws.Range("H3").Value = 70.488966
$ws.Range("I3").Value = 0.1321353991144917
$ws.Range("J3").Value = 0.1321353991144917
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.427350333333333
$ws.Range("N3").Value = 7.282051
$ws.Range("O3").Value = 0.2014397000898671
$ws.Range("P3").Value = 0.2014397000898671
$ws.Range("Q3").Value = 57.03380503880734
$ws.Range("R3").Value = 513.304245349266
$ws.Range("S3").Value = 0.02661731516887809
$ws.Range("T3").Value = 0.02661731516887809

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 23.496322
$ws.Range("H4").Value = 70.488966
$ws.Range("I4").Value = 0.1321353991144917
$ws.Range("J4").Value = 0.1321353991144917
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.785681666666667
$ws.Range("N4").Value = 26.357045
$ws.Range("O4").Value = 0.7291016281065776
$ws.Range("P4").Value = 0.7291016281065776
$ws.Range("Q4").Value = 206.4312054294967
$ws.Range("R4").Value = 1857.88084886547
$ws.Range("S4").Value = 0.0963401346248883
$ws.Range("T4").Value = 0.0963401346248883

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 124.9120333333333
$ws.Range("H5").Value = 374.7361
$ws.Range("I5").Value = 0.7024631931202969
$ws.Range("J5").Value = 0.7024631931202969
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.8369776666666665
$ws.Range("N5").Value = 2.510933
$ws.Range("O5").Value = 0.0694586718035551
$ws.Range("P5").Value = 0.06945867180355511
$ws.Range("Q5").Value = 104.5485821979222
$ws.Range("R5").Value = 940.9372397812998
$ws.Range("S5").Value = 0.04879216038502004
$ws.Range("T5").Value = 0.04879216038502006

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 124.9120333333333
$ws.Range("H6").Value = 374.7361
$ws.Range("I6").Value = 0.7024631931202969
$ws.Range("J6").Value = 0.7024631931202969
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.427350333333333
$ws.Range("N6").Value = 7.282051
$ws.Range("O6").Value = 0.2014397000898671
$ws.Range("P6").Value = 0.2014397000898671
$ws.Range("Q6").Value = 303.2052657490111
$ws.Range("R6").Value = 2728.8473917411
$ws.Range("S6").Value = 0.141503974946323
$ws.Range("T6").Value = 0.141503974946323

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 124.9120333333333
$ws.Range("H7").Value = 374.7361
$ws.Range("I7").Value = 0.7024631931202969
$ws.Range("J7").Value = 0.7024631931202969
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.785681666666667
$ws.Range("N7").Value = 26.357045
$ws.Range("O7").Value = 0.7291016281065776
$ws.Range("P7").Value = 0.7291016281065776
$ws.Range("Q7").Value = 1097.437361202722
$ws.Range("R7").Value = 9876.936250824499
$ws.Range("S7").Value = 0.5121670577889538
$ws.Range("T7").Value = 0.5121670577889538

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.411685
$ws.Range("H8").Value = 88.235055
$ws.Range("I8").Value = 0.1654014077652114
$ws.Range("J8").Value = 0.1654014077652114
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8369776666666665
$ws.Range("N8").Value = 2.510933
$ws.Range("O8").Value = 0.0694586718035551
$ws.Range("P8").Value = 0.06945867180355511
$ws.Range("Q8").Value = 24.616923484035
$ws.Range("R8").Value = 221.552311356315
$ws.Range("S8").Value = 0.01148856209780981
$ws.Range("T8").Value = 0.01148856209780981

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.411685
$ws.Range("H9").Value = 88.235055
$ws.Range("I9").Value = 0.1654014077652114
$ws.Range("J9").Value = 0.1654014077652114
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.427350333333333
$ws.Range("N9").Value = 7.282051
$ws.Range("O9").Value = 0.2014397000898671
$ws.Range("P9").Value = 0.2014397000898671
$ws.Range("Q9").Value = 71.39246338864501
$ws.Range("R9").Value = 642.5321704978051
$ws.Range("S9").Value = 0.03331840997466601
$ws.Range("T9").Value = 0.03331840997466601

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.411685
$ws.Range("H10").Value = 88.235055
$ws.Range("I10").Value = 0.1654014077652114
$ws.Range("J10").Value = 0.1654014077652114
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.785681666666667
$ws.Range("N10").Value = 26.357045
$ws.Range("O10").Value = 0.7291016281065776
$ws.Range("P10").Value = 0.7291016281065776
$ws.Range("Q10").Value = 258.401701690275
$ws.Range("R10").Value = 2325.615315212475
$ws.Range("S10").Value = 0.1205944356927356
$ws.Range("T10").Value = 0.1205944356927356
